$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 (Green LED, D1/D3): fill in Fabricant, Fournisseur, RefFabricant, RefFournisseur
# Leading "'" forces text entry (matches existing quote-prefixed cell style) without
# leaving a literal apostrophe in the stored value.
$ws.Range("D3").Value = "'KINGBRIGHT"
$ws.Range("F3").Value = "'Farnell"
$ws.Range("I3").Value = "'KPTD-2012LVZGCK"
$ws.Range("J3").Value = "'2846598"

# Row 4 (Red LED, D2): fill in Fabricant, Fournisseur, RefFabricant, RefFournisseur
$ws.Range("D4").Value = "'KINGBRIGHT"
$ws.Range("F4").Value = "'Farnell"
$ws.Range("I4").Value = "'KPTD-2012LVSURCK"
$ws.Range("J4").Value = "'2846595"
